$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-PlainValue($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

Set-PlainValue 2 4 "27.138.44"
Set-PlainValue 2 5 "  +0.37%  "
Set-PlainValue 3 4 "1.678.32"
Set-PlainValue 3 5 "  -0.12%  "
Set-PlainValue 4 5 "  +0.23%  "
Set-TextValue 5 4 "214.02"
Set-PlainValue 5 5 "  -0.91%  "
Set-PlainValue 6 5 "  -0.28%  "
Set-PlainValue 7 5 "  +0.22%  "
Set-PlainValue 8 5 "  +4.88%  "
Set-TextValue 9 4 "0.259"
Set-PlainValue 9 5 "  +1.87%  "
Set-TextValue 10 4 "0.0622"
Set-PlainValue 10 5 "  -0.29%  "
Set-TextValue 11 4 "0.0890"
Set-PlainValue 12 4 "1.918.35"
Set-PlainValue 12 5 "  +0.03%  "
Set-PlainValue 13 4 "1.686.70"
Set-PlainValue 13 5 "  +0.79%  "
Set-PlainValue 14 5 "  +1.72%  "
Set-TextValue 15 4 "0.554"
Set-PlainValue 15 5 "  +4.25%  "
Set-TextValue 16 4 "66.52"
Set-PlainValue 16 5 "  -0.03%  "
Set-PlainValue 17 4 "27.121.05"
Set-PlainValue 17 5 "  +0.32%  "
Set-TextValue 18 4 "235.27"
Set-PlainValue 18 5 "  -0.20%  "
Set-TextValue 19 4 "7.83"
Set-PlainValue 19 5 "  -4.25%  "
Set-PlainValue 20 5 "  -0.28%  "
Set-TextValue 22 4 "4.52"
Set-PlainValue 22 5 "  +1.24%  "
Set-TextValue 23 4 "9.51"
Set-PlainValue 23 5 "  +2.47%  "
Set-TextValue 24 4 "2.07"
Set-PlainValue 24 5 "  -2.20%  "
Set-TextValue 25 4 "146.62"
Set-PlainValue 25 5 "  -0.04%  "
Set-TextValue 26 4 "7.38"
Set-PlainValue 26 5 "  +1.86%  "
Set-TextValue 27 4 "16.30"
Set-PlainValue 28 5 "  -0.14%  "
Set-PlainValue 29 5 "  +0.24%  "
Set-TextValue 30 4 "0.0501"
Set-PlainValue 30 5 "  +0.60%  "
Set-PlainValue 31 5 "  -0.44%  "
Set-TextValue 32 4 "3.35"
Set-PlainValue 32 5 "  +0.08%  "
Set-PlainValue 33 4 "1.535.17"
Set-PlainValue 33 5 "  +0.60%  "
Set-PlainValue 34 5 "  +1.82%  "
Set-TextValue 35 4 "1.65"
Set-PlainValue 35 5 "  -3.52%  "
Set-TextValue 36 4 "0.602"
Set-PlainValue 36 5 "  +1.63%  "
Set-TextValue 37 4 "0.940"
Set-PlainValue 37 5 "  +2.14%  "
Set-TextValue 38 4 "2.40"
Set-PlainValue 38 5 "  -0.16%  "
Set-TextValue 39 4 "0.0171"
Set-PlainValue 39 5 "  -1.67%  "
Set-PlainValue 40 5 "  +3.88%  "
Set-PlainValue 41 5 "  +1.02%  "
Set-TextValue 42 4 "68.97"
Set-PlainValue 42 5 "  +1.63%  "
Set-PlainValue 43 5 "  +0.21%  "
Set-TextValue 44 4 "2.25"
Set-PlainValue 44 5 "  -0.09%  "
Set-PlainValue 45 4 "1.825.10"
Set-PlainValue 45 5 "  +0.21%  "
Set-TextValue 46 4 "0.790"
Set-PlainValue 46 5 "  +1.21%  "
Set-TextValue 47 4 "89.70"
Set-PlainValue 47 5 "  -0.61%  "
Set-PlainValue 48 5 "  +4.13%  "
Set-PlainValue 49 5 "  +6.28%  "
Set-TextValue 50 4 "8.18"
Set-TextValue 51 4 "0.103"
Set-PlainValue 51 5 "  -0.41%  "
